$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting rows 57..76 down to 58..77.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly price record.
$ws.Range("A57").Value = 11
$ws.Range("B57").Value = 'Vega Monumental Concepción'
$ws.Range("C57").Value = 'Bíobío'
$ws.Range("D57").Value = 44523
$ws.Range("E57").Value = 8
$ws.Range("F57").Value = 100112032
$ws.Range("G57").Value = 'Zapallo italiano'
$ws.Range("H57").Value = 'Sin especificar'
$ws.Range("I57").Value = 'Primera'
$ws.Range("J57").Value = 100
$ws.Range("K57").Value = 7500
$ws.Range("L57").Value = 8000
$ws.Range("M57").Value = 7750
$ws.Range("N57").Value = '$/caja 50 unidades'
$ws.Range("O57").Value = 'Región de Arica y Parinacota'
$ws.Range("P57").Value = 155
$ws.Range("Q57").Value = 50
$ws.Range("R57").Value = 'Hortaliza'
